# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list to use
# impact-focused accomplishment statements instead of job-duty style text,
# and trim the list down to four bullets (removing the ETL pipeline and
# AWS data warehouse duty-style bullets entirely).
#
# NOTE: several of these bullet sentences are duplicated verbatim elsewhere
# in the resume (e.g. under "PROFESSIONAL EXPERIENCE"), so a blind
# document-wide Find/Replace would also corrupt those untouched sections.
# Instead we walk the paragraph collection and only touch the exact
# paragraphs living in the "KEY ACHIEVEMENTS AND IMPACT" section.

$d = $word.ActiveDocument

# Map of old (unique, exact) paragraph text -> new text for the four
# bullets that get rewritten in place.
$replacements = @{
    "• Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%" = "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%";
    "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations" = "• `$4.7M savings enabled nonprofit access";
    "• Developed longitudinal data analysis methods using geospatial techniques that improved segmentation accuracy by 34% and survey incidence rates by 28%, reducing polling costs while increasing response quality" = "• 178% accuracy improvement in racial classification algorithms";
}

# Bullets to delete outright (whole paragraph incl. paragraph mark).
$deletions = @(
    "• Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets",
    "• Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy"
)

# Locate the "KEY ACHIEVEMENTS AND IMPACT" heading so we only operate on
# paragraphs that fall within that section (i.e. after this heading and
# before the next Heading2-level heading).
$paras = @($d.Paragraphs)
$sectionStartIndex = -1
for ($i = 0; $i -lt $paras.Count; $i++) {
    if ($paras[$i].Range.Text.Trim() -eq "KEY ACHIEVEMENTS AND IMPACT") {
        $sectionStartIndex = $i
        break
    }
}

if ($sectionStartIndex -ge 0) {
    $sectionEndIndex = $paras.Count - 1
    for ($i = $sectionStartIndex + 1; $i -lt $paras.Count; $i++) {
        if ($paras[$i].Style.NameLocal -eq "Heading 2") {
            $sectionEndIndex = $i - 1
            break
        }
    }

    # Handle the third bullet (the one duplicated elsewhere in the doc)
    # directly via its position: it is the paragraph immediately following
    # the "Built redistricting platform ... 89 organizations" bullet, which
    # is unique in the document.
    for ($i = $sectionStartIndex; $i -le $sectionEndIndex; $i++) {
        $txt = $paras[$i].Range.Text
        $trimmed = $txt.TrimEnd("`r", "`n", "`a")
        if ($trimmed -eq "• Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations") {
            $next = $paras[$i + 1]
            $nextTrimmed = $next.Range.Text.TrimEnd("`r", "`n", "`a")
            if ($nextTrimmed -eq "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis") {
                $rng = $next.Range
                $rng.Text = "• Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions"
            }
            break
        }
    }

    # Handle the remaining, uniquely-identifiable bullets with simple
    # whole-paragraph text replacement.
    foreach ($i in $sectionStartIndex..$sectionEndIndex) {
        $txt = $paras[$i].Range.Text
        $trimmed = $txt.TrimEnd("`r", "`n", "`a")
        if ($replacements.ContainsKey($trimmed)) {
            $paras[$i].Range.Text = $replacements[$trimmed]
        }
    }

    # Delete the two job-duty bullets entirely.
    foreach ($target in $deletions) {
        for ($i = $sectionEndIndex; $i -ge $sectionStartIndex; $i--) {
            $txt = $paras[$i].Range.Text
            $trimmed = $txt.TrimEnd("`r", "`n", "`a")
            if ($trimmed -eq $target) {
                $paras[$i].Range.Delete()
            }
        }
    }
}
